# Auto-generated Excel COM-interop script
# Applies the numeric corrections described by the commit diff
# across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2140.4
$ws.Range("J17").Value = 2140.4
$ws.Range("L17").Value = 6421.200000000001
$ws.Range("N17").Value = -6757.200000000001
$ws.Range("H64").Value = 4325.3335
$ws.Range("J64").Value = 3986.6667
$ws.Range("L64").Value = 3986.6667
$ws.Range("N64").Value = -4482.6667
$ws.Range("H67").Value = 4325.3335
$ws.Range("J67").Value = 3986.6667
$ws.Range("L67").Value = 3986.6667
$ws.Range("N67").Value = -5702.6667
$ws.Range("H100").Value = 2030.0588
$ws.Range("I100").Value = 2502.5
$ws.Range("K100").Value = 2502.5
$ws.Range("M100").Value = -1961.5
$ws.Range("H135").Value = 1051.7435
$ws.Range("I135").Value = 350.69696
$ws.Range("K135").Value = 3156.27264
$ws.Range("M135").Value = -621.2726400000001
$ws.Range("H137").Value = 1402.4828
$ws.Range("J137").Value = 1802.1333
$ws.Range("L137").Value = 5406.3999
$ws.Range("N137").Value = -10506.3999
$ws.Range("H138").Value = 525548.25
$ws.Range("I138").Value = 1737.6842
$ws.Range("J138").Value = 681054.5
$ws.Range("K138").Value = 5213.0526
$ws.Range("L138").Value = 2043163.5
$ws.Range("M138").Value = -73.05259999999998
$ws.Range("N138").Value = -2053443.5
$ws.Range("H141").Value = 626.4286
$ws.Range("I141").Value = 582.75
$ws.Range("K141").Value = 1748.25
$ws.Range("M141").Value = 3431.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3484.1125
$ws.Range("I32").Value = 3198.8333
$ws.Range("J32").Value = 6051.625
$ws.Range("K32").Value = 3198.8333
$ws.Range("L32").Value = 6051.625
$ws.Range("M32").Value = -2911.8333
$ws.Range("N32").Value = -6625.625
$ws.Range("H45").Value = 1274.52
$ws.Range("I45").Value = 1260.9333
$ws.Range("J45").Value = 1294.9
$ws.Range("K45").Value = 1260.9333
$ws.Range("L45").Value = 1294.9
$ws.Range("M45").Value = -883.9332999999999
$ws.Range("N45").Value = -2048.9
$ws.Range("H97").Value = 580.125
$ws.Range("I97").Value = 452.46155
$ws.Range("J97").Value = 1133.3334
$ws.Range("K97").Value = 452.46155
$ws.Range("L97").Value = 1133.3334
$ws.Range("M97").Value = 43.53845000000001
$ws.Range("N97").Value = -2125.3334
$ws.Range("H102").Value = 33334198.0
$ws.Range("I102").Value = 33334198.0
$ws.Range("K102").Value = 33334198.0
$ws.Range("M102").Value = -33332576.0
$ws.Range("H134").Value = 34999.5
$ws.Range("J134").Value = 34999.5
$ws.Range("L134").Value = 34999.5
$ws.Range("N134").Value = -45139.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 12500751.0
$ws.Range("I94").Value = 15625645.0
$ws.Range("K94").Value = 15625645.0
$ws.Range("M94").Value = -15625194.0
$ws.Range("H105").Value = 252473970.0
$ws.Range("I105").Value = 336631300.0
$ws.Range("J105").Value = 2000.0
$ws.Range("K105").Value = 336631300.0
$ws.Range("L105").Value = 2000.0
$ws.Range("M105").Value = -336629553.0
$ws.Range("N105").Value = -5494.0
$ws.Range("H134").Value = 7154.636
$ws.Range("J134").Value = 26400.0
$ws.Range("L134").Value = 79200.0
$ws.Range("N134").Value = -84270.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1491.1111
$ws.Range("I31").Value = 1282.0
$ws.Range("K31").Value = 1282.0
$ws.Range("M31").Value = -987.0
$ws.Range("H34").Value = 1491.1111
$ws.Range("I34").Value = 1282.0
$ws.Range("K34").Value = 1282.0
$ws.Range("M34").Value = -1080.0
$ws.Range("H86").Value = 3936868.0
$ws.Range("I86").Value = 7410454.5
$ws.Range("K86").Value = 7410454.5
$ws.Range("M86").Value = -7409331.5
$ws.Range("H89").Value = 3936868.0
$ws.Range("I89").Value = 7410454.5
$ws.Range("K89").Value = 37052272.5
$ws.Range("M89").Value = -37046656.5
$ws.Range("H129").Value = 37999.4
$ws.Range("I129").Value = 0.0
$ws.Range("J129").Value = 37999.4
$ws.Range("K129").Value = 0.0
$ws.Range("L129").Value = 37999.4
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -47999.4
$ws.Range("H132").Value = 1752.9445
$ws.Range("I132").Value = 1387.0
$ws.Range("J132").Value = 3582.6667
$ws.Range("K132").Value = 4161.0
$ws.Range("L132").Value = 10748.0001
$ws.Range("M132").Value = -1631.0
$ws.Range("N132").Value = -15808.0001
$ws.Range("H134").Value = 1159.2142
$ws.Range("I134").Value = 1041.375
$ws.Range("J134").Value = 1316.3334
$ws.Range("K134").Value = 3124.125
$ws.Range("L134").Value = 3949.0002
$ws.Range("M134").Value = -589.125
$ws.Range("N134").Value = -9019.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2463.6924
$ws.Range("J81").Value = 3064.375
$ws.Range("L81").Value = 9193.125
$ws.Range("N81").Value = -11439.125
$ws.Range("H84").Value = 2463.6924
$ws.Range("J84").Value = 3064.375
$ws.Range("L84").Value = 27579.375
$ws.Range("N84").Value = -38811.375
$ws.Range("H137").Value = 23442806.0
$ws.Range("I137").Value = 41668030.0
$ws.Range("J137").Value = 10374.5
$ws.Range("K137").Value = 125004090.0
$ws.Range("L137").Value = 31123.5
$ws.Range("M137").Value = -124998990.0
$ws.Range("N137").Value = -41323.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H129").Value = 47499.25
$ws.Range("I129").Value = 40000.0
$ws.Range("J129").Value = 49999.0
$ws.Range("K129").Value = 40000.0
$ws.Range("L129").Value = 49999.0
$ws.Range("M129").Value = -35000.0
$ws.Range("N129").Value = -59999.0
$ws.Range("H132").Value = 1820.0256
$ws.Range("I132").Value = 1198.6923
$ws.Range("K132").Value = 3596.0769
$ws.Range("M132").Value = -1066.0769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 920.2222
$ws.Range("I100").Value = 683.0
$ws.Range("K100").Value = 683.0
$ws.Range("M100").Value = -142.0
$ws.Range("H124").Value = 35000.0
$ws.Range("J124").Value = 35000.0
$ws.Range("L124").Value = 35000.0
$ws.Range("N124").Value = -44820.0
$ws.Range("H132").Value = 19113.578
$ws.Range("I132").Value = 1024.9143
$ws.Range("K132").Value = 3074.7429
$ws.Range("M132").Value = -544.7428999999997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2318.5366
$ws.Range("I132").Value = 2341.1667
$ws.Range("J132").Value = 2155.6
$ws.Range("K132").Value = 7023.500100000001
$ws.Range("L132").Value = 6466.799999999999
$ws.Range("M132").Value = -4493.500100000001
$ws.Range("N132").Value = -11526.8
$ws.Range("H136").Value = 409.94736
$ws.Range("I136").Value = 222.82353
$ws.Range("K136").Value = 668.47059
$ws.Range("M136").Value = 1881.52941
